# Add Asian-handicap ("亚盘") odds columns alongside the existing European odds
# export. New columns AA:AG are appended after the existing A:Z data on the
# product sheet, with corresponding header labels added to row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new headers (columns AA:AG) ---
$ws.Cells.Item(1, 27).Value = "亚盘公司"
$ws.Cells.Item(1, 28).Value = "初盘主队水位"
$ws.Cells.Item(1, 29).Value = "初盘盘口"
$ws.Cells.Item(1, 30).Value = "初盘客队水位"
$ws.Cells.Item(1, 31).Value = "即时主队水位"
$ws.Cells.Item(1, 32).Value = "即时盘口"
$ws.Cells.Item(1, 33).Value = "即时客队水位"

# --- Row 2 data ---
$ws.Cells.Item(2, 27).Value = "澳门"
$ws.Cells.Item(2, 28).Value = "0.84"
$ws.Cells.Item(2, 29).Value = "一球/球半"
$ws.Cells.Item(2, 30).Value = "1.02"
$ws.Cells.Item(2, 31).Value = "0.76"
$ws.Cells.Item(2, 32).Value = "一球"
$ws.Cells.Item(2, 33).Value = "1.10"

# --- Row 3 data ---
$ws.Cells.Item(3, 27).Value = "澳门"
$ws.Cells.Item(3, 28).Value = "0.84"
$ws.Cells.Item(3, 29).Value = "一球/球半"
$ws.Cells.Item(3, 30).Value = "1.02"
$ws.Cells.Item(3, 31).Value = "0.76"
$ws.Cells.Item(3, 32).Value = "一球"
$ws.Cells.Item(3, 33).Value = "1.10"

# --- Row 4 data ---
$ws.Cells.Item(4, 27).Value = "澳门"
$ws.Cells.Item(4, 28).Value = "0.84"
$ws.Cells.Item(4, 29).Value = "一球/球半"
$ws.Cells.Item(4, 30).Value = "1.02"
$ws.Cells.Item(4, 31).Value = "0.76"
$ws.Cells.Item(4, 32).Value = "一球"
$ws.Cells.Item(4, 33).Value = "1.10"
